# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$xlHAlignCenter = -4108

# 1. Insert "Player Info" sheet as the first sheet
$firstSheet = $wb.Worksheets.Item(1)
$playerInfoSheet = $wb.Worksheets.Add($firstSheet)
$playerInfoSheet.Name = "Player Info"

$playerInfoData = @(
    @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE"),
    @("3991", "Bhuvneshwar Kumar", "Right Handed", "Right Arm Medium")
)

for ($r = 0; $r -lt $playerInfoData.Length; $r++) {
    $rowData = $playerInfoData[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $cell = $playerInfoSheet.Cells.Item($r + 1, $c + 1)
            $cell.Value = $val
            if ($r -eq 0) {
                $cell.Font.Bold = $true
                $cell.HorizontalAlignment = $xlHAlignCenter
                $cell.VerticalAlignment = -4160
                $cell.Borders.LineStyle = 1
            }
        }
    }
}

# 2. Rewrite "ODI Batting" sheet contents
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Clear()

$battingData = @(
    @("MATCH_NUMBER", "INNING_NUMBER", "MATCH_DATE", "MATCH_CODE", "MATCH_INNING", "OPPONENT", "VENUE", "DISMISSAL", "RUNS_SCORED", "BALLS_FACED"),
    @("1", $null, "30/12/2012", "3455", "1st", "Pakistan", "MA Chidambaram Stadium", "did not bat", "-", "-"),
    @("2", "1", "03/01/2013", "3456", "2nd", "Pakistan", "Eden Gardens", "lbw b Saeed Ajmal", "0", "2"),
    @("3", "2", "06/01/2013", "3457", "1st", "Pakistan", "Arun Jaitley Stadium", "lbw b Saeed Ajmal", "2", "10"),
    @("4", "3", "11/01/2013", "3459", "2nd", "England", "Saurashtra Cricket Association Stadium", "not out", "20*", "16"),
    @("5", $null, "15/01/2013", "3461", "1st", "England", "Nehru Stadium (Kochi)", "did not bat", "-", "-"),
    @("6", $null, "19/01/2013", "3463", "2nd", "England", "JSCA International Stadium Complex", "did not bat", "-", "-"),
    @("7", $null, "23/01/2013", "3468", "2nd", "England", "Punjab Cricket Association IS Bindra Stadium", "did not bat", "-", "-"),
    @("8", "4", "27/01/2013", "3470", "1st", "England", "Himachal Pradesh Cricket Association Stadium", "c S T Finn b T T Bresnan", "31", "30"),
    @("9", "5", "06/06/2013", "3505", "1st", "South Africa", "Sophia Gardens", "not out", "0*", "0"),
    @("10", $null, "11/06/2013", "3510", "2nd", "West Indies", "Kennington Oval", "did not bat", "-", "-"),
    @("11", $null, "15/06/2013", "3514", "2nd", "Pakistan", "Edgbaston", "did not bat", "-", "-"),
    @("12", $null, "20/06/2013", "3518", "2nd", "Sri Lanka", "Sophia Gardens", "did not bat", "-", "-"),
    @("13", "6", "23/06/2013", "3519", "1st", "England", "Edgbaston", "not out", "1*", "1"),
    @("14", "7", "30/06/2013", "3521", "1st", "West Indies", "Sabina Park", "not out", "11*", "5"),
    @("15", $null, "05/07/2013", "3525", "1st", "West Indies", "Queen's Park Oval", "did not bat", "-", "-"),
    @("16", $null, "09/07/2013", "3529", "1st", "Sri Lanka", "Queen's Park Oval", "did not bat", "-", "-"),
    @("17", "8", "11/07/2013", "3530", "2nd", "Sri Lanka", "Queen's Park Oval", "lbw b S L Malinga", "0", "15"),
    @("18", "9", "13/10/2013", "3562", "2nd", "Australia", "Maharashtra Cricket Association Stadium", "c A C Voges b A J Finch", "18", "27"),
    @("19", $null, "16/10/2013", "3563", "2nd", "Australia", "Sawai Mansingh Stadium", "did not bat", "-", "-"),
    @("20", "10", "19/10/2013", "3564", "1st", "Australia", "Punjab Cricket Association IS Bindra Stadium", "c G J Bailey b J P Faulkner", "10", "18"),
    @("21", $null, "30/10/2013", "3568", "2nd", "Australia", "Vidarbha Cricket Association Stadium", "did not bat", "-", "-"),
    @("22", $null, "02/11/2013", "3572", "1st", "Australia", "M Chinnaswamy Stadium", "did not bat", "-", "-"),
    @("23", $null, "21/11/2013", "3580", "2nd", "West Indies", "Nehru Stadium (Kochi)", "did not bat", "-", "-"),
    @("24", "11", "24/11/2013", "3581", "1st", "West Indies", "Dr YS Rajasekhara Reddy Cricket Stadium", "not out", "1*", "1"),
    @("25", $null, "27/11/2013", "3583", "2nd", "West Indies", "Green Park", "did not bat", "-", "-"),
    @("26", "12", "05/12/2013", "3586", "2nd", "South Africa", "Wanderers Stadium", "c J H Kallis b D W Steyn", "0", "4"),
    @("27", "13", "19/01/2014", "3601", "2nd", "New Zealand", "McLean Park", "run out", "6", "4"),
    @("28", "14", "22/01/2014", "3603", "2nd", "New Zealand", "Seddon Park", "c N L McCullum b C J Anderson", "11", "6"),
    @("29", "15", "25/01/2014", "3607", "2nd", "New Zealand", "Eden Park", "c N L McCullum b H K Bennett", "4", "3"),
    @("30", $null, "28/01/2014", "3609", "1st", "New Zealand", "Seddon Park", "did not bat", "-", "-"),
    @("31", "16", "31/01/2014", "3612", "2nd", "New Zealand", "Sky Stadium", "c †L Ronchi b M J Henry", "20", "25"),
    @("32", $null, "26/02/2014", "3619", "2nd", "Bangladesh", "Khan Shaheb Osman Ali Stadium ", "did not bat", "-", "-"),
    @("33", "17", "28/02/2014", "3621", "1st", "Sri Lanka", "Khan Shaheb Osman Ali Stadium ", "st K C Sangakkara b B A W Mendis", "0", "0"),
    @("34", $null, "01/03/2014", "3624", "1st", "Pakistan", "Shere Bangla National Stadium", "did not bat", "-", "-"),
    @("35", $null, "05/03/2014", "3628", "2nd", "Afghanistan", "Shere Bangla National Stadium", "did not bat", "-", "-"),
    @("36", $null, "27/08/2014", "3664", "1st", "England", "Sophia Gardens", "did not bat", "-", "-"),
    @("37", $null, "30/08/2014", "3667", "2nd", "England", "Trent Bridge", "did not bat", "-", "-"),
    @("38", $null, "02/09/2014", "3670", "2nd", "England", "Edgbaston", "did not bat", "-", "-"),
    @("39", "18", "05/09/2014", "3672", "2nd", "England", "Headingley", "run out", "1", "3"),
    @("40", "19", "08/10/2014", "3678", "2nd", "West Indies", "Nehru Stadium (Kochi)", "c D J G Sammy b M N Samuels", "2", "10"),
    @("41", "20", "11/10/2014", "3680", "1st", "West Indies", "Arun Jaitley Stadium", "c K A Pollard b D J Bravo", "18", "14"),
    @("42", $null, "17/10/2014", "3683", "1st", "West Indies", "Himachal Pradesh Cricket Association Stadium", "did not bat", "-", "-"),
    @("43", "21", "18/01/2015", "3731", "1st", "Australia", "Melbourne Cricket Ground", "b M A Starc", "0", "1"),
    @("44", "22", "20/01/2015", "3735", "1st", "England", "Brisbane Cricket Ground", "b J M Anderson", "5", "5"),
    @("45", $null, "28/02/2015", "3768", "2nd", "United Arab Emirates", "WACA Ground", "did not bat", "-", "-"),
    @("46", "23", "18/06/2015", "3808", "2nd", "Bangladesh", "Shere Bangla National Stadium", "not out", "25*", "21"),
    @("47", "24", "21/06/2015", "3810", "1st", "Bangladesh", "Shere Bangla National Stadium", "c †Litton Das b Rubel Hossain", "3", "10"),
    @("48", $null, "10/07/2015", "3812", "1st", "Zimbabwe", "Harare Sports Club", "did not bat", "-", "-"),
    @("49", "25", "12/07/2015", "3815", "1st", "Zimbabwe", "Harare Sports Club", "not out", "0*", "0"),
    @("50", $null, "14/07/2015", "3817", "1st", "Zimbabwe", "Harare Sports Club", "did not bat", "-", "-"),
    @("51", "26", "11/10/2015", "3841", "2nd", "South Africa", "Green Park", "not out", "1*", "1"),
    @("52", "27", "14/10/2015", "3844", "1st", "South Africa", "Holkar Cricket Stadium", "b Imran Tahir", "14", "32"),
    @("53", $null, "18/10/2015", "3845", "2nd", "South Africa", "Saurashtra Cricket Association Stadium", "did not bat", "-", "-"),
    @("54", "28", "22/10/2015", "3848", "1st", "South Africa", "MA Chidambaram Stadium", "run out", "0", "1"),
    @("55", "29", "25/10/2015", "3851", "2nd", "South Africa", "Wankhede Stadium", "c D A Miller b Imran Tahir", "1", "8"),
    @("56", $null, "12/01/2016", "3874", "1st", "Australia", "WACA Ground", "did not bat", "-", "-"),
    @("57", "30", "20/01/2016", "3877", "2nd", "Australia", "Manuka Oval", "c S P D Smith b K W Richardson", "2", "6"),
    @("58", $null, "19/01/2017", "3976", "1st", "England", "Barabati Stadium", "did not bat", "-", "-"),
    @("59", "31", "22/01/2017", "3978", "2nd", "England", "Eden Gardens", "not out", "0*", "1"),
    @("60", $null, "04/06/2017", "4034", "1st", "Pakistan", "Edgbaston", "did not bat", "-", "-"),
    @("61", $null, "08/06/2017", "4038", "1st", "Sri Lanka", "Kennington Oval", "did not bat", "-", "-"),
    @("62", $null, "11/06/2017", "4042", "2nd", "South Africa", "Kennington Oval", "did not bat", "-", "-"),
    @("63", $null, "15/06/2017", "4047", "2nd", "Bangladesh", "Edgbaston", "did not bat", "-", "-"),
    @("64", "32", "18/06/2017", "4050", "2nd", "Pakistan", "Kennington Oval", "not out", "1*", "8"),
    @("65", $null, "23/06/2017", "4051", "1st", "West Indies", "Queen's Park Oval", "did not bat", "-", "-"),
    @("66", $null, "25/06/2017", "4052", "1st", "West Indies", "Queen's Park Oval", "did not bat", "-", "-"),
    @("67", $null, "30/06/2017", "4053", "1st", "West Indies", "Sir Vivian Richards Stadium", "did not bat", "-", "-"),
    @("68", $null, "20/08/2017", "4061", "2nd", "Sri Lanka", "Rangiri Dambulla International Stadium", "did not bat", "-", "-"),
    @("69", "33", "24/08/2017", "4062", "2nd", "Sri Lanka", "Pallekele International Cricket Stadium", "not out", "53*", "80"),
    @("70", $null, "27/08/2017", "4063", "2nd", "Sri Lanka", "Pallekele International Cricket Stadium", "did not bat", "-", "-"),
    @("71", $null, "03/09/2017", "4065", "2nd", "Sri Lanka", "R Premadasa Stadium", "did not bat", "-", "-"),
    @("72", "34", "17/09/2017", "4067", "1st", "Australia", "MA Chidambaram Stadium", "not out", "32*", "30"),
    @("73", "35", "21/09/2017", "4069", "1st", "Australia", "Eden Gardens", "c G J Maxwell b K W Richardson", "20", "33"),
    @("74", $null, "24/09/2017", "4071", "2nd", "Australia", "Holkar Cricket Stadium", "did not bat", "-", "-"),
    @("75", $null, "01/10/2017", "4076", "2nd", "Australia", "Vidarbha Cricket Association Stadium", "did not bat", "-", "-"),
    @("76", "36", "22/10/2017", "4085", "1st", "New Zealand", "Wankhede Stadium", "c H M Nicholls b T G Southee", "26", "15"),
    @("77", $null, "25/10/2017", "4088", "2nd", "New Zealand", "Maharashtra Cricket Association Stadium", "did not bat", "-", "-"),
    @("78", $null, "29/10/2017", "4089", "1st", "New Zealand", "Green Park", "did not bat", "-", "-"),
    @("79", "37", "10/12/2017", "4096", "1st", "Sri Lanka", "Himachal Pradesh Cricket Association Stadium", "c †N Dickwella b R A S Lakmal", "0", "5"),
    @("80", $null, "13/12/2017", "4098", "1st", "Sri Lanka", "Punjab Cricket Association IS Bindra Stadium", "did not bat", "-", "-"),
    @("81", $null, "17/12/2017", "4099", "2nd", "Sri Lanka", "Dr YS Rajasekhara Reddy Cricket Stadium", "did not bat", "-", "-"),
    @("82", $null, "01/02/2018", "4126", "2nd", "South Africa", "Kingsmead", "did not bat", "-", "-"),
    @("83", $null, "04/02/2018", "4127", "2nd", "South Africa", "SuperSport Park", "did not bat", "-", "-"),
    @("84", "38", "07/02/2018", "4128", "1st", "South Africa", "Newlands", "not out", "16*", "19"),
    @("85", "39", "10/02/2018", "4130", "1st", "South Africa", "Wanderers Stadium", "run out", "5", "7"),
    @("86", "40", "13/02/2018", "4133", "1st", "South Africa", "St George's Park", "not out", "19*", "20"),
    @("87", "41", "17/07/2018", "4175", "1st", "England", "Headingley", "c J M Bairstow b D J Willey", "21", "35"),
    @("88", "42", "18/09/2018", "4196", "1st", "Hong Kong", "Dubai International Cricket Stadium", "c Anshuman Rath b K D Shah", "9", "18"),
    @("89", $null, "19/09/2018", "4197", "2nd", "Pakistan", "Dubai International Cricket Stadium", "did not bat", "-", "-"),
    @("90", $null, "21/09/2018", "4199", "2nd", "Bangladesh", "Dubai International Cricket Stadium", "did not bat", "-", "-"),
    @("91", $null, "23/09/2018", "4201", "2nd", "Pakistan", "Dubai International Cricket Stadium", "did not bat", "-", "-"),
    @("92", "43", "28/09/2018", "4205", "2nd", "Bangladesh", "Dubai International Cricket Stadium", "c †Mushfiqur Rahim b Mustafizur Rahman", "21", "31"),
    @("93", "44", "27/10/2018", "4219", "2nd", "West Indies", "Maharashtra Cricket Association Stadium", "c R Powell b O C McCoy", "10", "17"),
    @("94", $null, "29/10/2018", "4220", "1st", "West Indies", "Brabourne Stadium", "did not bat", "-", "-"),
    @("95", $null, "01/11/2018", "4221", "2nd", "West Indies", "Greenfield International Stadium", "did not bat", "-", "-"),
    @("96", "45", "12/01/2019", "4234", "2nd", "Australia", "Sydney Cricket Ground", "not out", "29*", "23"),
    @("97", $null, "15/01/2019", "4235", "2nd", "Australia", "Adelaide Oval", "did not bat", "-", "-"),
    @("98", $null, "18/01/2019", "4236", "2nd", "Australia", "Melbourne Cricket Ground", "did not bat", "-", "-"),
    @("99", $null, "23/01/2019", "4239", "2nd", "New Zealand", "McLean Park", "did not bat", "-", "-"),
    @("100", $null, "26/01/2019", "4242", "1st", "New Zealand", "Bay Oval", "did not bat", "-", "-"),
    @("101", $null, "28/01/2019", "4245", "2nd", "New Zealand", "Bay Oval", "did not bat", "-", "-"),
    @("102", "46", "31/01/2019", "4248", "1st", "New Zealand", "Seddon Park", "b C de Grandhomme", "1", "12"),
    @("103", "47", "03/02/2019", "4249", "1st", "New Zealand", "Sky Stadium", "c L R P L Taylor b T A Boult", "6", "8"),
    @("104", "48", "10/03/2019", "4268", "1st", "Australia", "Punjab Cricket Association IS Bindra Stadium", "c †A T Carey b J A Richardson", "1", "2"),
    @("105", "49", "13/03/2019", "4270", "2nd", "Australia", "Arun Jaitley Stadium", "c A J Finch b P J Cummins", "46", "54"),
    @("106", $null, "05/06/2019", "4310", "2nd", "South Africa", "The Rose Bowl", "did not bat", "-", "-"),
    @("107", $null, "09/06/2019", "4316", "1st", "Australia", "Kennington Oval", "did not bat", "-", "-"),
    @("108", $null, "16/06/2019", "4324", "1st", "Pakistan", "Old Trafford", "did not bat", "-", "-"),
    @("109", "50", "02/07/2019", "4345", "1st", "Bangladesh", "Edgbaston", "run out", "2", "3"),
    @("110", $null, "06/07/2019", "4350", "2nd", "Sri Lanka", "Headingley", "did not bat", "-", "-"),
    @("111", "51", "09/07/2019", "4353", "2nd", "New Zealand", "Old Trafford", "b L H Ferguson", "0", "1"),
    @("112", $null, "08/08/2019", "4359", " ", "West Indies", "Providence Stadium", "did not bat", "-", "-"),
    @("113", "52", "11/08/2019", "4360", "1st", "West Indies", "Queen's Park Oval", "c K A J Roach b C R Brathwaite", "1", "2"),
    @("114", $null, "14/08/2019", "4362", "2nd", "West Indies", "Queen's Park Oval", "did not bat", "-", "-"),
    @("115", $null, "23/03/2021", "4454", "1st", "England", "Maharashtra Cricket Association Stadium", "did not bat", "-", "-"),
    @("116", $null, "26/03/2021", "4456", "1st", "England", "Maharashtra Cricket Association Stadium", "did not bat", "-", "-"),
    @("117", "53", "28/03/2021", "4457", "1st", "England", "Maharashtra Cricket Association Stadium", "c S M Curran b R J W Topley", "3", "5"),
    @("118", $null, "18/07/2021", "4480", "2nd", "Sri Lanka", "R Premadasa Stadium", "did not bat", "-", "-"),
    @("119", "54", "20/07/2021", "4482", "2nd", "Sri Lanka", "R Premadasa Stadium", "not out", "19*", "28"),
    @("120", "55", "19/01/2022", "4524", "2nd", "South Africa", "Boland Park", "c T Bavuma b T Shamsi", "4", "11"),
    @("121", $null, "21/01/2022", "4526", "1st", "South Africa", "Boland Park", "did not bat", "-", "-")
)

for ($r = 0; $r -lt $battingData.Length; $r++) {
    $rowData = $battingData[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $cell = $battingSheet.Cells.Item($r + 1, $c + 1)
            $cell.Value = $val
            if ($r -eq 0) {
                $cell.Font.Bold = $true
                $cell.HorizontalAlignment = $xlHAlignCenter
                $cell.VerticalAlignment = -4160
                $cell.Borders.LineStyle = 1
            }
        }
    }
}

# 3. Rewrite "ODI Bowling" sheet contents
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Clear()

$bowlingData = @(
    @("MATCH_NUMBER", "MATCH_CODE", "MATCH_INNING", "OPPONENT", "VENUE", "OVERS", "WICKETS_RUNS"),
    @("1", "3455", "2nd", "Pakistan", "MA Chidambaram Stadium", "9.0", "2/27"),
    @("2", "3456", "1st", "Pakistan", "Eden Gardens", "9.0", "1/61"),
    @("3", "3457", "2nd", "Pakistan", "Arun Jaitley Stadium", "10.0", "2/31"),
    @("4", "3459", "1st", "England", "Saurashtra Cricket Association Stadium", "7.0", "0/52"),
    @("5", "3461", "2nd", "England", "Nehru Stadium (Kochi)", "10.0", "3/29"),
    @("6", "3463", "1st", "England", "JSCA International Stadium Complex", "10.0", "1/40"),
    @("7", "3468", "1st", "England", "Punjab Cricket Association IS Bindra Stadium", "10.0", "0/30"),
    @("8", "3470", "2nd", "England", "Himachal Pradesh Cricket Association Stadium", "9.0", "0/45"),
    @("9", "3505", "2nd", "South Africa", "Sophia Gardens", "7.0", "2/49"),
    @("10", "3510", "1st", "West Indies", "Kennington Oval", "8.0", "1/32"),
    @("11", "3514", "1st", "Pakistan", "Edgbaston", "8.0", "2/19"),
    @("12", "3518", "1st", "Sri Lanka", "Sophia Gardens", "9.0", "1/18"),
    @("13", "3519", "2nd", "England", "Edgbaston", "3.0", "0/19"),
    @("14", "3521", "2nd", "West Indies", "Sabina Park", "7.0", "1/36"),
    @("15", "3525", "2nd", "West Indies", "Queen's Park Oval", "8.0", "3/29"),
    @("16", "3529", "2nd", "Sri Lanka", "Queen's Park Oval", "6.0", "4/8"),
    @("17", "3530", "1st", "Sri Lanka", "Queen's Park Oval", "8.0", "2/24"),
    @("18", "3562", "1st", "Australia", "Maharashtra Cricket Association Stadium", "7.0", "0/41"),
    @("19", "3563", "1st", "Australia", "Sawai Mansingh Stadium", "10.0", "0/54"),
    @("20", "3564", "2nd", "Australia", "Punjab Cricket Association IS Bindra Stadium", "10.0", "1/50"),
    @("21", "3568", "1st", "Australia", "Vidarbha Cricket Association Stadium", "8.0", "1/42"),
    @("22", "3572", "2nd", "Australia", "M Chinnaswamy Stadium", "8.0", "0/47"),
    @("23", "3580", "1st", "West Indies", "Nehru Stadium (Kochi)", "5.0", "0/26"),
    @("24", "3581", "2nd", "West Indies", "Dr YS Rajasekhara Reddy Cricket Stadium", "9.0", "2/56"),
    @("25", "3583", "1st", "West Indies", "Green Park", "8.0", "1/42"),
    @("26", "3586", "1st", "South Africa", "Wanderers Stadium", "9.0", "0/68"),
    @("27", "3601", "1st", "New Zealand", "McLean Park", "10.0", "1/38"),
    @("28", "3603", "1st", "New Zealand", "Seddon Park", "7.0", "1/43"),
    @("29", "3607", "1st", "New Zealand", "Eden Park", "9.0", "1/48"),
    @("30", "3609", "2nd", "New Zealand", "Seddon Park", "10.0", "0/62"),
    @("31", "3612", "1st", "New Zealand", "Sky Stadium", "8.0", "1/48"),
    @("32", "3619", "1st", "Bangladesh", "Khan Shaheb Osman Ali Stadium ", "8.0", "1/41"),
    @("33", "3621", "2nd", "Sri Lanka", "Khan Shaheb Osman Ali Stadium ", "9.2", "0/45"),
    @("34", "3624", "2nd", "Pakistan", "Shere Bangla National Stadium", "10.0", "2/56"),
    @("35", "3628", "1st", "Afghanistan", "Shere Bangla National Stadium", "8.0", "0/25"),
    @("36", "3664", "2nd", "England", "Sophia Gardens", "7.0", "1/30"),
    @("37", "3667", "1st", "England", "Trent Bridge", "8.0", "1/45"),
    @("38", "3670", "1st", "England", "Edgbaston", "8.0", "2/14"),
    @("39", "3672", "1st", "England", "Headingley", "8.0", "1/45"),
    @("40", "3678", "1st", "West Indies", "Nehru Stadium (Kochi)", "10.0", "0/38"),
    @("41", "3680", "2nd", "West Indies", "Arun Jaitley Stadium", "7.0", "0/32"),
    @("42", "3683", "2nd", "West Indies", "Himachal Pradesh Cricket Association Stadium", "10.0", "2/25"),
    @("43", "3731", "2nd", "Australia", "Melbourne Cricket Ground", "9.5", "1/44"),
    @("44", "3735", "2nd", "England", "Brisbane Cricket Ground", "2.0", "0/18"),
    @("45", "3768", "1st", "United Arab Emirates", "WACA Ground", "5.0", "1/19"),
    @("46", "3808", "1st", "Bangladesh", "Shere Bangla National Stadium", "7.0", "2/37"),
    @("47", "3810", "2nd", "Bangladesh", "Shere Bangla National Stadium", "5.0", "0/32"),
    @("48", "3812", "2nd", "Zimbabwe", "Harare Sports Club", "10.0", "1/35"),
    @("49", "3815", "2nd", "Zimbabwe", "Harare Sports Club", "10.0", "4/33"),
    @("50", "3817", "2nd", "Zimbabwe", "Harare Sports Club", "6.0", "0/12"),
    @("51", "3841", "1st", "South Africa", "Green Park", "10.0", "0/67"),
    @("52", "3844", "2nd", "South Africa", "Holkar Cricket Stadium", "8.4", "3/41"),
    @("53", "3845", "1st", "South Africa", "Saurashtra Cricket Association Stadium", "10.0", "0/65"),
    @("54", "3848", "2nd", "South Africa", "MA Chidambaram Stadium", "10.0", "3/68"),
    @("55", "3851", "1st", "South Africa", "Wankhede Stadium", "10.0", "1/106"),
    @("56", "3874", "2nd", "Australia", "WACA Ground", "9.0", "0/42"),
    @("57", "3877", "1st", "Australia", "Manuka Oval", "8.0", "0/69"),
    @("58", "3976", "2nd", "England", "Barabati Stadium", "10.0", "1/63"),
    @("59", "3978", "1st", "England", "Eden Gardens", "8.0", "0/56"),
    @("60", "4034", "2nd", "Pakistan", "Edgbaston", "5.0", "1/23"),
    @("61", "4038", "2nd", "Sri Lanka", "Kennington Oval", "10.0", "1/54"),
    @("62", "4042", "1st", "South Africa", "Kennington Oval", "7.3", "2/23"),
    @("63", "4047", "1st", "Bangladesh", "Edgbaston", "10.0", "2/53"),
    @("64", "4050", "1st", "Pakistan", "Kennington Oval", "10.0", "1/44"),
    @("66", "4052", "2nd", "West Indies", "Queen's Park Oval", "5.0", "2/9"),
    @("67", "4053", "2nd", "West Indies", "Sir Vivian Richards Stadium", "5.0", "0/19"),
    @("68", "4061", "1st", "Sri Lanka", "Rangiri Dambulla International Stadium", "6.0", "0/33"),
    @("69", "4062", "1st", "Sri Lanka", "Pallekele International Cricket Stadium", "10.0", "0/53"),
    @("70", "4063", "1st", "Sri Lanka", "Pallekele International Cricket Stadium", "9.0", "0/41"),
    @("71", "4065", "1st", "Sri Lanka", "R Premadasa Stadium", "9.4", "5/42"),
    @("72", "4067", "2nd", "Australia", "MA Chidambaram Stadium", "4.0", "1/25"),
    @("73", "4069", "2nd", "Australia", "Eden Gardens", "6.1", "3/9"),
    @("74", "4071", "1st", "Australia", "Holkar Cricket Stadium", "10.0", "0/52"),
    @("75", "4076", "1st", "Australia", "Vidarbha Cricket Association Stadium", "8.0", "1/40"),
    @("76", "4085", "2nd", "New Zealand", "Wankhede Stadium", "10.0", "1/56"),
    @("77", "4088", "1st", "New Zealand", "Maharashtra Cricket Association Stadium", "10.0", "3/45"),
    @("78", "4089", "2nd", "New Zealand", "Green Park", "10.0", "1/92"),
    @("79", "4096", "2nd", "Sri Lanka", "Himachal Pradesh Cricket Association Stadium", "8.4", "1/42"),
    @("80", "4098", "2nd", "Sri Lanka", "Punjab Cricket Association IS Bindra Stadium", "9.0", "1/40"),
    @("81", "4099", "1st", "Sri Lanka", "Dr YS Rajasekhara Reddy Cricket Stadium", "6.5", "1/35"),
    @("82", "4126", "1st", "South Africa", "Kingsmead", "10.0", "1/71"),
    @("83", "4127", "1st", "South Africa", "SuperSport Park", "5.0", "1/19"),
    @("84", "4128", "2nd", "South Africa", "Newlands", "7.0", "0/41"),
    @("85", "4130", "2nd", "South Africa", "Wanderers Stadium", "4.0", "0/27"),
    @("86", "4133", "2nd", "South Africa", "St George's Park", "7.0", "0/43"),
    @("87", "4175", "2nd", "England", "Headingley", "7.0", "0/49"),
    @("88", "4196", "2nd", "Hong Kong", "Dubai International Cricket Stadium", "9.0", "0/50"),
    @("89", "4197", "1st", "Pakistan", "Dubai International Cricket Stadium", "7.0", "3/15"),
    @("90", "4199", "1st", "Bangladesh", "Dubai International Cricket Stadium", "10.0", "3/32"),
    @("91", "4201", "1st", "Pakistan", "Dubai International Cricket Stadium", "9.0", "0/46"),
    @("92", "4205", "1st", "Bangladesh", "Dubai International Cricket Stadium", "7.0", "0/33"),
    @("93", "4219", "1st", "West Indies", "Maharashtra Cricket Association Stadium", "10.0", "1/70"),
    @("94", "4220", "2nd", "West Indies", "Brabourne Stadium", "5.0", "1/30"),
    @("95", "4221", "1st", "West Indies", "Greenfield International Stadium", "4.0", "1/11"),
    @("96", "4234", "1st", "Australia", "Sydney Cricket Ground", "10.0", "2/66"),
    @("97", "4235", "1st", "Australia", "Adelaide Oval", "10.0", "4/45"),
    @("98", "4236", "1st", "Australia", "Melbourne Cricket Ground", "8.0", "2/28"),
    @("99", "4239", "1st", "New Zealand", "McLean Park", "5.0", "0/20"),
    @("100", "4242", "2nd", "New Zealand", "Bay Oval", "7.0", "2/42"),
    @("101", "4245", "1st", "New Zealand", "Bay Oval", "10.0", "2/46"),
    @("102", "4248", "2nd", "New Zealand", "Seddon Park", "5.0", "2/25"),
    @("103", "4249", "2nd", "New Zealand", "Sky Stadium", "7.1", "1/38"),
    @("104", "4268", "2nd", "Australia", "Punjab Cricket Association IS Bindra Stadium", "9.0", "1/67"),
    @("105", "4270", "1st", "Australia", "Arun Jaitley Stadium", "10.0", "3/48"),
    @("106", "4310", "1st", "South Africa", "The Rose Bowl", "10.0", "2/44"),
    @("107", "4316", "2nd", "Australia", "Kennington Oval", "10.0", "3/50"),
    @("108", "4324", "2nd", "Pakistan", "Old Trafford", "2.4", "0/8"),
    @("109", "4345", "2nd", "Bangladesh", "Edgbaston", "9.0", "1/51"),
    @("110", "4350", "1st", "Sri Lanka", "Headingley", "10.0", "1/73"),
    @("111", "4353", "1st", "New Zealand", "Old Trafford", "10.0", "3/43"),
    @("112", "4359", "1st", "West Indies", "Providence Stadium", "4.0", "0/11"),
    @("113", "4360", "2nd", "West Indies", "Queen's Park Oval", "8.0", "4/31"),
    @("114", "4362", "1st", "West Indies", "Queen's Park Oval", "5.0", "0/48"),
    @("115", "4454", "2nd", "England", "Maharashtra Cricket Association Stadium", "9.0", "2/30"),
    @("116", "4456", "2nd", "England", "Maharashtra Cricket Association Stadium", "10.0", "1/63"),
    @("117", "4457", "2nd", "England", "Maharashtra Cricket Association Stadium", "10.0", "3/42"),
    @("118", "4480", "1st", "Sri Lanka", "R Premadasa Stadium", "9.0", "0/63"),
    @("119", "4482", "1st", "Sri Lanka", "R Premadasa Stadium", "10.0", "3/54"),
    @("120", "4524", "1st", "South Africa", "Boland Park", "10.0", "0/64"),
    @("121", "4526", "2nd", "South Africa", "Boland Park", "8.0", "0/67")
)

for ($r = 0; $r -lt $bowlingData.Length; $r++) {
    $rowData = $bowlingData[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $cell = $bowlingSheet.Cells.Item($r + 1, $c + 1)
            $cell.Value = $val
            if ($r -eq 0) {
                $cell.Font.Bold = $true
                $cell.HorizontalAlignment = $xlHAlignCenter
                $cell.VerticalAlignment = -4160
                $cell.Borders.LineStyle = 1
            }
        }
    }
}

# 4. Append "ODI Batting Extra" sheet as the last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

$extraData = @(
    @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH"),
    @("4248", 8, "0", "0", "1.09%", "NO"),
    @("4249", 9, "0", "0", "2.38%", "NO"),
    @("4268", 8, "0", "0", "0.28%", "NO"),
    @("4270", 8, "3", "2", "19.41%", "NO"),
    @("4310", $null, $null, $null, $null, "NO"),
    @("4316", 8, $null, $null, $null, "NO"),
    @("4324", 8, $null, $null, $null, "NO"),
    @("4345", 8, "0", "0", "0.64%", "NO"),
    @("4350", 9, $null, $null, $null, "NO"),
    @("4353", $null, $null, $null, $null, "NO"),
    @("4359", $null, $null, $null, $null, "NO"),
    @("4360", 8, "0", "0", "0.36%", "NO"),
    @("4362", $null, $null, $null, $null, "NO"),
    @("4454", 9, $null, $null, $null, "NO"),
    @("4456", $null, $null, $null, $null, "NO"),
    @("4457", 9, "0", "0", "0.91%", "NO"),
    @("4480", 8, $null, $null, $null, "NO"),
    @("4482", 9, "2", "0", "6.86%", "NO"),
    @("4524", $null, $null, $null, $null, "NO"),
    @("4526", 9, $null, $null, $null, "NO")
)

for ($r = 0; $r -lt $extraData.Length; $r++) {
    $rowData = $extraData[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($val -ne $null) {
            $cell = $extraSheet.Cells.Item($r + 1, $c + 1)
            $cell.Value = $val
            if ($r -eq 0) {
                $cell.Font.Bold = $true
                $cell.HorizontalAlignment = $xlHAlignCenter
                $cell.VerticalAlignment = -4160
                $cell.Borders.LineStyle = 1
            }
        }
    }
}

Write-Host "Edit complete. Sheets:" ($wb.Worksheets | ForEach-Object { $_.Name })
